# Weekly update: a new price observation is inserted at row 475, pushing the
# existing rows 475-504 down to 476-505 (dimension grows from A1:T504 to
# A1:T505). The new row 475 carries the latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 475, shifting rows 475:504
# down to 476:505.
$ws.Rows.Item(475).Insert()

# Populate the newly inserted row 475 with the new weekly observation.
$ws.Range("A475").Value = 2
$ws.Range("B475").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C475").Value = "Coquimbo"
$ws.Range("D475").Value = 44931
$ws.Range("E475").Value = 4
$ws.Range("F475").Value = "Fruta"
$ws.Range("G475").Value = 100102
$ws.Range("H475").Value = "Cítricos"
$ws.Range("I475").Value = 100102005
$ws.Range("J475").Value = "Naranja"
$ws.Range("K475").Value = "Valencia"
$ws.Range("L475").Value = "Primera"
$ws.Range("M475").Value = 16
$ws.Range("N475").Value = 200000
$ws.Range("O475").Value = 210000
$ws.Range("P475").Value = 205000
$ws.Range("Q475").Value = "`$/bins (400 kilos)"
$ws.Range("R475").Value = "Provincia de Limarí"
$ws.Range("S475").Value = 512
$ws.Range("T475").Value = 400
